# Update the "F" column (想去人数 / interest count) figures that changed
# between data pulls on sheets "展览" and "全部类型".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 203
$ws1.Range("F6").Value  = 9724
$ws1.Range("F10").Value = 2886
$ws1.Range("F11").Value = 167
$ws1.Range("F12").Value = 109
$ws1.Range("F13").Value = 24
$ws1.Range("F15").Value = 278
$ws1.Range("F16").Value = 508
$ws1.Range("F18").Value = 263
$ws1.Range("F19").Value = 1398

# Sheet "全部类型" (all types) - same events, rows offset by +1
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 203
$ws4.Range("F7").Value  = 9724
$ws4.Range("F11").Value = 2886
$ws4.Range("F12").Value = 167
$ws4.Range("F13").Value = 109
$ws4.Range("F14").Value = 24
$ws4.Range("F16").Value = 278
$ws4.Range("F17").Value = 508
$ws4.Range("F19").Value = 263
$ws4.Range("F20").Value = 1398
